$d = $word.ActiveDocument

$replacements = @(
    @{old="75×37="; new="31×33="},
    @{old="60×81="; new="67×39="},
    @{old="38×66="; new="67×91="},
    @{old="61×95="; new="43×99="},
    @{old="46×40="; new="56×40="},
    @{old="69×17="; new="17×61="},
    @{old="68×13="; new="68×82="},
    @{old="55×51="; new="79×18="},
    @{old="39×61="; new="60×30="},
    @{old="63×33="; new="82×51="},
    @{old="96×37="; new="59×44="},
    @{old="24×99="; new="48×72="},
    @{old="91×36="; new="32×97="},
    @{old="88×90="; new="81×79="},
    @{old="48×50="; new="84×38="},
    @{old="90×50="; new="97×53="},
    @{old="33×35="; new="70×26="},
    @{old="43×71="; new="36×24="},
    @{old="11×12="; new="87×45="},
    @{old="61×31="; new="48×84="},
    @{old="34×63="; new="91×41="},
    @{old="50×52="; new="45×27="},
    @{old="51×65="; new="49×33="},
    @{old="20×49="; new="14×26="},
    @{old="43×65="; new="18×86="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
